$d = $word.ActiveDocument

# The document contains five numbered ("Paragrafoelenco") list items followed
# by one trailing empty list paragraph. The commit removes the four
# "to-do" bullet items that already got actioned (footer/form fixes) plus
# the "whatsapp"-related one, leaving only a single empty paragraph that is
# no longer part of the numbered list (but keeps the "Paragrafoelenco"
# paragraph style).

# Delete the four populated list paragraphs (repeatedly removing paragraph 1
# shifts the remaining ones up).
$d.Paragraphs.Item(1).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()
$d.Paragraphs.Item(1).Range.Delete()

# Only the originally-last (empty) paragraph remains now. Strip its list
# numbering (<w:numPr>) while keeping the "Paragrafoelenco" paragraph style,
# matching the target markup. Switching the style away and back drops the
# paragraph's direct numPr formatting.
$d.Paragraphs.Item(1).Style = "Normal"
$d.Paragraphs.Item(1).Style = "Paragrafoelenco"
